$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (AFONSO BINGALA JONE) ---
$ws.Range("D2").Value = "PARQUE DE GORONGOSA"
$ws.Range("E2").Value = "TECNICO DE POUPANCA"
$ws.Range("F2").Value = "'"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'  9/8/2022"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value = 1
$ws.Range("K2").Value = "'"
$ws.Range("K2").Style = "Normal"

# --- Row 3 (ANTONIO AGOSTINHO JOAO NOBRE) ---
$ws.Range("D3").Value = "PARQUE DE GORONGOSA"
$ws.Range("E3").Value = "TECNICO DE CAFE"
$ws.Range("F3").Value = "'"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "'  9/8/2022"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value = 1
$ws.Range("K3").Value = "'"
$ws.Range("K3").Style = "Normal"

# --- Row 4 (FLORINDA NETO) ---
$ws.Range("D4").Value = "GORONGOSA"
$ws.Range("E4").Value = "TENICO"
$ws.Range("F4").Value = "'847233663"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Value = "'10/22/2022"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value = 4

# --- Row 5 (JOSSEFO CELESTINO SALIVA) ---
$ws.Range("D5").Value = "GORONGOSA"
$ws.Range("E5").Value = "TECNICO DO CAMPO"
$ws.Range("F5").Value = "'869210890"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Value = "'10/21/2022"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value = 3

# --- Row 6 (JULEIDA ZULFA CARLOS) ---
$ws.Range("D6").Value = "GORONGOSA"
$ws.Range("E6").Value = "TECNICA"
$ws.Range("F6").Value = "'863019855"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Value = "'10/22/2022"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value = 4
$ws.Range("K6").Value = "'84038571"
$ws.Range("K6").Style = "Normal"
